$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("K2").Value = 10
$ws.Range("N2").Value = 2.03
$ws.Range("O2").Value = 1.78

# Row 6 updates
$ws.Range("G6").Value = 1.52
$ws.Range("H6").Value = 4.8
$ws.Range("I6").Value = 4.7
$ws.Range("L6").Value = 1.07
$ws.Range("M6").Value = 6.8
$ws.Range("N6").Value = 1.25
$ws.Range("O6").Value = 3.55
$ws.Range("R6").Value = 1.35
$ws.Range("S6").Value = 2.95
$ws.Range("U6").Value = 13.5
$ws.Range("W6").Value = 15
$ws.Range("X6").Value = 10.75
$ws.Range("Y6").Value = 14.5
$ws.Range("Z6").Value = 40
$ws.Range("AA6").Value = 12
$ws.Range("AB6").Value = 13
$ws.Range("AC6").Value = 28
$ws.Range("AD6").Value = 110
$ws.Range("AE6").Value = 32
$ws.Range("AF6").Value = 45
$ws.Range("AG6").Value = 17
$ws.Range("AH6").Value = 90
